# Adds the 2021年 data row to Sheet1, mirroring the existing
# 2018年/2019年/2020年 rows (A2:DK4) already present in the workbook.
# Columns B:DK hold the year-over-year growth figures for each industry.
# Column AG has no reported figure for any year in this sheet, so (like
# AG2:AG4) it is written as an explicit empty-text cell rather than a
# numeric 0 or a cell that is simply left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row goes right after the last row currently used on the sheet
# (rows 1-4 -> header + 2018/2019/2020 -> new row is 5).
$newRow = $ws.UsedRange.Rows.Count + 1
$lastCol = 115   # column DK
$yearColNum = 33 # column AG

$rowValues = New-Object 'object[,]' 1,115
$rowValues[0,0] = "2021年"
$rowValues[0,1] = 7.5
$rowValues[0,2] = 23.7
$rowValues[0,3] = -21.4
$rowValues[0,4] = 52.4
$rowValues[0,5] = 0.6
$rowValues[0,6] = 81.2
$rowValues[0,7] = 14.3
$rowValues[0,8] = 6.5
$rowValues[0,9] = 7
$rowValues[0,10] = -8.9
$rowValues[0,11] = -12.2
$rowValues[0,12] = 13.8
$rowValues[0,13] = -38.2
$rowValues[0,14] = -1.3
$rowValues[0,15] = 6
$rowValues[0,16] = -10.2
$rowValues[0,17] = 30.7
$rowValues[0,18] = 15.4
$rowValues[0,19] = 10.5
$rowValues[0,20] = 10.2
$rowValues[0,21] = 4.7
$rowValues[0,22] = 18.2
$rowValues[0,23] = 14.6
$rowValues[0,24] = 15.2
$rowValues[0,25] = 41.8
$rowValues[0,26] = 10
$rowValues[0,27] = 25.3
$rowValues[0,28] = 20.1
$rowValues[0,29] = 7.8
$rowValues[0,30] = 12.6
$rowValues[0,31] = -40.8
$rowValues[0,32] = $null
$rowValues[0,33] = -23.5
$rowValues[0,34] = -8.300000000000001
$rowValues[0,35] = -16.9
$rowValues[0,36] = 16
$rowValues[0,37] = 4.9
$rowValues[0,38] = 0.2
$rowValues[0,39] = -9.5
$rowValues[0,40] = -10.3
$rowValues[0,41] = -25.3
$rowValues[0,42] = 7
$rowValues[0,43] = 1.3
$rowValues[0,44] = 43.7
$rowValues[0,45] = 52.7
$rowValues[0,46] = 45.8
$rowValues[0,47] = 4.6
$rowValues[0,48] = 11
$rowValues[0,49] = -6.8
$rowValues[0,50] = -7.2
$rowValues[0,51] = 11.4
$rowValues[0,52] = 1.8
$rowValues[0,53] = 6.1
$rowValues[0,54] = 8.6
$rowValues[0,55] = 58.5
$rowValues[0,56] = 4.6
$rowValues[0,57] = 3
$rowValues[0,58] = 13.6
$rowValues[0,59] = -3.3
$rowValues[0,60] = -3.9
$rowValues[0,61] = 13.5
$rowValues[0,62] = 13.3
$rowValues[0,63] = -1.2
$rowValues[0,64] = 1.2
$rowValues[0,65] = -4.6
$rowValues[0,66] = 1.3
$rowValues[0,67] = 4.2
$rowValues[0,68] = 35.3
$rowValues[0,69] = 10.7
$rowValues[0,70] = -4.4
$rowValues[0,71] = -2.2
$rowValues[0,72] = -10.5
$rowValues[0,73] = 1.7
$rowValues[0,74] = 4.5
$rowValues[0,75] = 22.8
$rowValues[0,76] = 20.2
$rowValues[0,77] = 1
$rowValues[0,78] = 7.3
$rowValues[0,79] = -9.4
$rowValues[0,80] = 16.4
$rowValues[0,81] = -28.2
$rowValues[0,82] = -3.3
$rowValues[0,83] = 14.7
$rowValues[0,84] = 17.4
$rowValues[0,85] = 6.2
$rowValues[0,86] = 12.3
$rowValues[0,87] = -15.5
$rowValues[0,88] = 11.6
$rowValues[0,89] = 5.4
$rowValues[0,90] = -18
$rowValues[0,91] = 18.8
$rowValues[0,92] = 12.7
$rowValues[0,93] = 26.7
$rowValues[0,94] = 10.3
$rowValues[0,95] = -25
$rowValues[0,96] = 20.7
$rowValues[0,97] = 9.9
$rowValues[0,98] = 11.8
$rowValues[0,99] = -1.3
$rowValues[0,100] = 19
$rowValues[0,101] = 18.9
$rowValues[0,102] = 7.3
$rowValues[0,103] = 34.7
$rowValues[0,104] = 11.4
$rowValues[0,105] = 1.8
$rowValues[0,106] = 19.4
$rowValues[0,107] = -2.1
$rowValues[0,108] = -7.5
$rowValues[0,109] = 13.5
$rowValues[0,110] = 26.4
$rowValues[0,111] = 9.300000000000001
$rowValues[0,112] = 10.5
$rowValues[0,113] = 15.4
$rowValues[0,114] = 27.7

$startCell = $ws.Cells.Item($newRow, 1)
$endCell = $ws.Cells.Item($newRow, $lastCol)
$ws.Range($startCell, $endCell).Value = $rowValues

# AG column (e.g. AG5) must exist as an empty *text* cell (matches
# AG2:AG4), not just be skipped by the array write above ($null values
# are left untouched by Range.Value). A leading apostrophe forces Excel
# to store an explicit empty string there instead of leaving the cell
# absent; then reset the quote-prefix formatting it implies so the cell
# keeps the sheet's default (unstyled) look, same as the other AG cells.
$yearCell = $ws.Cells.Item($newRow, $yearColNum)
$yearCell.Value = "'"
$yearCell.Style = "Normal"

# Give the new row's label cell (A5) the same look as the existing
# year labels in A2:A4 (bold, centered, bordered).
$labelCell = $ws.Cells.Item($newRow, 1)
$priorLabelCell = $ws.Cells.Item($newRow - 1, 1)
$priorLabelCell.Copy() | Out-Null
$labelCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

Write-Host ("Added new year row {0} through column {1}." -f $newRow, $lastCol)
